# Uppercase the airport codes in the "From" (F) and "To" (G) columns so
# values like "Blr"/"Del" read as "BLR"/"DEL" throughout the flights table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $fVal = $fCell.Value2
    $gVal = $gCell.Value2

    if ($fVal -ne $null) {
        $fCell.Value2 = $fVal.ToUpper()
    }
    if ($gVal -ne $null) {
        $gCell.Value2 = $gVal.ToUpper()
    }
}
